# Append: 2025-11-10 12:50 JST
# Update the "取得日時" (retrieved datetime) column (A) for all existing
# data rows on the "ランサーズ" sheet from the previous run timestamp to
# the new one, reflecting a refreshed scrape/append pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-10 12:47:45"
$newTimestamp = "2025-11-10 12:50:05"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
